$d = $word.ActiveDocument

# 1. Split the first paragraph's text ("{#content}{text}") into three runs
#    with proofErr gramStart/gramEnd markers bracketing the middle run,
#    matching what Word's grammar checker inserts when that span is re-edited.
$p1 = $d.Paragraphs.Item(1).Range
$p1.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>{#</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>content}{</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>text}</w:t></w:r></w:p>')

# 2. Remove the blank paragraph, the "Repository: <hyperlink>" paragraph and
#    the trailing blank paragraph that followed "{/content}".
$start = $d.Paragraphs.Item(3).Range.Start
$end = $d.Paragraphs.Item(5).Range.End
$d.Range($start, $end).Delete()
